# Add a "Save" column (H) to the s_vals sheet, matching the style used
# by the other header/value columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1 header: copy G1's formatting (bold font, border, center/top align)
# so the new header cell reuses the same style instead of creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# H2:H8 data values from the diff
$saveValues = @(0, 1, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
